$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "raju"
$ws.Range("B3").Value = "khanna"

$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:B3"))

$ws.Range("A3").Select()
